# Update the furniture template to the new ASSET style:
#  - add two new columns to the header row: furniture_type (M) and radius (N)
#  - mark the first TABLE (row 5) as a "rectangle" furniture_type with rotation 0
#  - mark the second TABLE (row 6) as a "circle" furniture_type with a radius,
#    removing its now-unused width/height values
#  - leave the active selection on L14, matching the author's last click

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns (order matters: it controls shared-string allocation
# order so it matches the authored workbook - furniture_type, rectangle,
# circle, then radius).
$ws.Range("M1").Value = "furniture_type"

# Row 5 (TABLE @ 26x34.5): rectangle-style furniture, explicit rotation.
$ws.Range("G5").Value = 0
$ws.Range("M5").Value = "rectangle"

# Row 6 (TABLE @ 30x80): circle-style furniture - drop width/height, add radius.
$ws.Range("E6:F6").ClearContents()
$ws.Range("M6").Value = "circle"
$ws.Range("N6").Value = 10

# Finish the header row with the radius column.
$ws.Range("N1").Value = "radius"

# Restore the author's last selection.
$ws.Range("L14").Select()
